$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add P1 and Q1 headers, continuing the 0..15 sequence, with same style as O1 (header style)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("P1:Q1").Font.Bold = $true
$ws.Range("P1:Q1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("P1:Q1").VerticalAlignment = -4160     # xlTop
$ws.Range("P1:Q1").Borders.LineStyle = 1         # xlContinuous

# For each data row (2..25): set new columns P,Q = 2, and swap I<->K and M<->O values
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I : was 1 -> 2
    $ws.Cells.Item($r, 11).Value = 1  # K : was 2 -> 1
    $ws.Cells.Item($r, 13).Value = 2  # M : was 1 -> 2
    $ws.Cells.Item($r, 15).Value = 1  # O : was 2 -> 1
    $ws.Cells.Item($r, 16).Value = 2  # P : new
    $ws.Cells.Item($r, 17).Value = 2  # Q : new
}
